$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.079.11'
$ws.Range("E2").Value = '  +0.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.890.22'
$ws.Range("E3").Value = '  +1.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5155'
$ws.Range("E7").Value = '  +2.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3755'
$ws.Range("E8").Value = '  +3.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07196'
$ws.Range("E9").Value = '  +0.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.14'
$ws.Range("E10").Value = '  +1.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9012'
$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("E12").Value = '  +2.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.880.44'
$ws.Range("E13").Value = '  +1.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.25'
$ws.Range("E14").Value = '  +1.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.247'
$ws.Range("E15").Value = '  +0.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008496'
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.47'
$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.125.71'
$ws.Range("E20").Value = '  +0.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.061'
$ws.Range("E21").Value = '  +0.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.106.81'
$ws.Range("E22").Value = '  +1.28%  '

$ws.Range("E23").Value = '  +1.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.379'
$ws.Range("E24").Value = '  -0.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.289'
$ws.Range("E25").Value = '  +10.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.52'
$ws.Range("E26").Value = '  -0.40%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.10'
$ws.Range("E27").Value = '  +1.20%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.732'
$ws.Range("E28").Value = '  -3.23%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.37'
$ws.Range("E29").Value = '  +1.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.921'
$ws.Range("E30").Value = '  +5.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.798'
$ws.Range("E31").Value = '  +2.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09208'
$ws.Range("E32").Value = '  -0.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05056'
$ws.Range("E33").Value = '  -0.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.235'
$ws.Range("E34").Value = '  +7.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7692'
$ws.Range("E35").Value = '  +2.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.985'
$ws.Range("E36").Value = '  -0.25%  '

$ws.Range("E37").Value = '  +0.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.601'
$ws.Range("E38").Value = '  +2.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5614'
$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01994'
$ws.Range("E40").Value = '  -0.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.070'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.115'
$ws.Range("E42").Value = '  +6.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.641'
$ws.Range("E43").Value = '  +1.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.52'
$ws.Range("E44").Value = '  -0.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1504'
$ws.Range("E45").Value = '  +2.28%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4822'
$ws.Range("E46").Value = '  +2.66%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.20'
$ws.Range("E47").Value = '  +1.17%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9997'
$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.609'
$ws.Range("E49").Value = '  +2.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.63'
$ws.Range("E50").Value = '  +2.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.02'
$ws.Range("E51").Value = '  +1.75%  '
